$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update score values in column B
$ws.Range("B3").Value = 86
$ws.Range("B4").Value = 86
$ws.Range("B7").Value = 75

# Update the active selection to B7
$ws.Range("B7").Select()
